# Updated cryptos list (prices + 1h volume %) per GitHub Actions scrape run.
# Rows 24/25 swap ranking: Stellar <-> Toncoin trade places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "19.56"); Excel's COM
# layer auto-coerces such literals to floating-point numbers on a bare
# `.Value =` assignment (losing the exact text + introducing FP noise like
# 215.86000000000001). Force text type via NumberFormat, assign, then restore
# the cell style to Normal so no stray formatting sticks around.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.000.85'
$ws.Range("E2").Value = '  +0.63%  '
Set-TextValue $ws.Range("D3") '1.641.70'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E4").Value = '  +0.39%  '
Set-TextValue $ws.Range("D5") '215.86'
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("E8").Value = '  +0.50%  '
$ws.Range("E9").Value = '  +1.07%  '
Set-TextValue $ws.Range("D10") '19.56'
$ws.Range("E10").Value = '  +0.17%  '
Set-TextValue $ws.Range("D11") '0.0797'
$ws.Range("E11").Value = '  +0.69%  '
Set-TextValue $ws.Range("D12") '1.869.19'
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("E13").Value = '  +0.54%  '
Set-TextValue $ws.Range("D14") '1.646.86'
$ws.Range("E14").Value = '  -1.37%  '
Set-TextValue $ws.Range("D15") '0.544'
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("E16").Value = '  +1.09%  '
Set-TextValue $ws.Range("D17") '63.37'
$ws.Range("E17").Value = '  +1.27%  '
Set-TextValue $ws.Range("D18") '26.102.77'
$ws.Range("E18").Value = '  +1.00%  '
$ws.Range("E19").Value = '  +0.40%  '
Set-TextValue $ws.Range("D20") '194.40'
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("E21").Value = '  -0.62%  '
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("E23").Value = '  -0.75%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D24") '1.80'
$ws.Range("E24").Value = '  -0.98%  '
$ws.Range("B25").Value = 'Stellar'
$ws.Range("C25").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D25") '0.131'
$ws.Range("E25").Value = '  +4.64%  '
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("E27").Value = '  -0.20%  '
Set-TextValue $ws.Range("D28") '6.88'
$ws.Range("E28").Value = '  +0.76%  '
Set-TextValue $ws.Range("D29") '15.52'
$ws.Range("E29").Value = '  +0.75%  '
$ws.Range("E30").Value = '  +0.77%  '
Set-TextValue $ws.Range("D31") '0.0495'
$ws.Range("E31").Value = '  -0.39%  '
$ws.Range("E32").Value = '  -0.15%  '
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("E34").Value = '  -1.04%  '
$ws.Range("E35").Value = '  +1.17%  '
$ws.Range("E36").Value = '  +0.48%  '
Set-TextValue $ws.Range("D37") '1.129.58'
$ws.Range("E37").Value = '  -0.75%  '
Set-TextValue $ws.Range("D38") '0.539'
$ws.Range("E38").Value = '  -1.32%  '
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("E40").Value = '  +0.29%  '
Set-TextValue $ws.Range("D41") '5.48'
$ws.Range("E41").Value = '  +0.84%  '
Set-TextValue $ws.Range("D42") '99.05'
$ws.Range("E42").Value = '  -0.04%  '
Set-TextValue $ws.Range("D43") '0.797'
$ws.Range("E43").Value = '  +0.20%  '
Set-TextValue $ws.Range("D44") '1.777.84'
$ws.Range("E44").Value = '  +0.69%  '
$ws.Range("E45").Value = '  +4.02%  '
Set-TextValue $ws.Range("D46") '56.43'
$ws.Range("E46").Value = '  +0.46%  '
Set-TextValue $ws.Range("D47") '0.0522'
$ws.Range("E47").Value = '  -1.09%  '
$ws.Range("E48").Value = '  +3.52%  '
Set-TextValue $ws.Range("D49") '7.72'
$ws.Range("E49").Value = '  +0.98%  '
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("E51").Value = '  -0.45%  '
